# Sample Project / Main.xlsx - "Rules" sheet
# Row 11 (the R40 rule row) - change the Rule name from "R40" to "1".
# The leading apostrophe forces Excel to keep the numeric-looking entry
# stored as text (matching the original cell's text type) instead of
# silently converting it to the number 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("B11").Value = "'1"
